$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "fixed bug for swap scene": the CanClone flag (column K) was stuck at 0
# for every scene row except the template row (10). Turn it on for every
# real scene (rows 11-35). Row 12 (DemoSummer) had the Share flag
# (column J) incorrectly left on, so it is swapped off as part of the fix.

$ws.Cells.Item(12, 10).Value = 0   # J12 (Share) : 1 -> 0

for ($row = 11; $row -le 35; $row++) {
    $ws.Cells.Item($row, 11).Value = 1   # column K (CanClone) : 0 -> 1
}

# Restore the user's on-screen selection after the edit.
$ws.Range("K13").Select()
